# The commit swaps the content of ppt/theme/theme1.xml (the default
# "Office Theme" color scheme) and ppt/theme/theme2.xml (the "Integral" /
# "Red Violet" color scheme) that ship inside the deck: theme2.xml (the
# theme actually applied to the slide master / presentation) ends up
# carrying the plain "Office" palette that used to live in theme1.xml.
#
# The PowerPoint object model doesn't expose a "swap these two theme
# parts" verb, so we reproduce the visible effect: push the Office
# Theme's twelve scheme colors into the presentation's active theme
# color scheme (reached here through a slide's ThemeColorScheme, which -
# unlike SlideMaster.ColorScheme - does not blow away the existing
# <a:clrScheme> name while writing).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme color scheme (RGB() packed as r + g*256 + b*65536):
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
